# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# for a handful of Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, matching freshly pulled Universalis pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9566
$ws.Range("L43").Value = 11041.889
$ws.Range("N43").Value = -11179.889
$ws.Range("M43").Value = -2855.5
$ws.Range("I43").Value = 2924.5
$ws.Range("J43").Value = 11041.889
$ws.Range("K43").Value = 2924.5

$ws.Range("H74").Value = 3664.5
$ws.Range("L74").Value = 3996.75
$ws.Range("N74").Value = -5868.75
$ws.Range("M74").Value = -2064
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 3996.75
$ws.Range("K74").Value = 3000

$ws.Range("H77").Value = 3664.5
$ws.Range("L77").Value = 19983.75
$ws.Range("N77").Value = -29343.75
$ws.Range("M77").Value = -10320
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 3996.75
$ws.Range("K77").Value = 15000

$ws.Range("H80").Value = 126004.19
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7996
$ws.Range("M80").Value = -401815.42
$ws.Range("I80").Value = 134271.14
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 402813.42

$ws.Range("H83").Value = 126004.19
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27984
$ws.Range("M83").Value = -1203448.26
$ws.Range("I83").Value = 134271.14
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 1208440.26

$ws.Range("H112").Value = 2317.238
$ws.Range("L112").Value = 7438.5
$ws.Range("N112").Value = -9654.5
$ws.Range("J112").Value = 2479.5

$ws.Range("H115").Value = 867.1818
$ws.Range("M115").Value = 200.5
$ws.Range("I115").Value = 455.5
$ws.Range("K115").Value = 1366.5

$ws.Range("H132").Value = 11354.818
$ws.Range("L132").Value = 18300
$ws.Range("N132").Value = -23360
$ws.Range("M132").Value = -44671.499
$ws.Range("I132").Value = 15733.833
$ws.Range("J132").Value = 6100
$ws.Range("K132").Value = 47201.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6298
$ws.Range("L61").Value = 5998
$ws.Range("N61").Value = -6422
$ws.Range("M61").Value = -6186
$ws.Range("I61").Value = 6398
$ws.Range("J61").Value = 5998
$ws.Range("K61").Value = 6398

$ws.Range("H135").Value = 349258.66
$ws.Range("L135").Value = 349258.66
$ws.Range("N135").Value = -359398.66
$ws.Range("J135").Value = 349258.66

$ws.Range("H136").Value = 6298
$ws.Range("L136").Value = 17994
$ws.Range("N136").Value = -23094
$ws.Range("M136").Value = -16644
$ws.Range("I136").Value = 6398
$ws.Range("J136").Value = 5998
$ws.Range("K136").Value = 19194

$ws.Range("H139").Value = 139501
$ws.Range("L139").Value = 139501
$ws.Range("N139").Value = -149781
$ws.Range("J139").Value = 139501

$ws.Range("H141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360
$ws.Range("J141").Value = 150000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 6000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -5832
$ws.Range("I12").Value = 6000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 6000
$ws.Range("N12").ClearContents()

$ws.Range("H99").Value = 20670
$ws.Range("M99").Value = -29886.084
$ws.Range("I99").Value = 31384.084
$ws.Range("K99").Value = 31384.084

$ws.Range("H105").Value = 76432.21000000001
$ws.Range("L105").Value = 7039.4287
$ws.Range("N105").Value = -10533.4287
$ws.Range("M105").Value = -144078
$ws.Range("I105").Value = 145825
$ws.Range("J105").Value = 7039.4287
$ws.Range("K105").Value = 145825

$ws.Range("H138").Value = 96420.86
$ws.Range("L138").Value = 96420.86
$ws.Range("N138").Value = -106700.86
$ws.Range("J138").Value = 96420.86

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11617116
$ws.Range("L99").Value = 7899.25
$ws.Range("N99").Value = -10895.25
$ws.Range("M99").Value = -19355096
$ws.Range("I99").Value = 19356594
$ws.Range("J99").Value = 7899.25
$ws.Range("K99").Value = 19356594

$ws.Range("H126").Value = 11617116
$ws.Range("L126").Value = 23697.75
$ws.Range("N126").Value = -28637.75
$ws.Range("M126").Value = -58067312
$ws.Range("I126").Value = 19356594
$ws.Range("J126").Value = 7899.25
$ws.Range("K126").Value = 58069782

$ws.Range("H134").Value = 2185.4
$ws.Range("L134").Value = 5098.7142
$ws.Range("N134").Value = -10168.7142
$ws.Range("M134").Value = -4587.999899999999
$ws.Range("I134").Value = 2374.3333
$ws.Range("J134").Value = 1699.5714
$ws.Range("K134").Value = 7122.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 48972740
$ws.Range("L4").Value = 316280400
$ws.Range("N4").Value = -316280624
$ws.Range("J4").Value = 105426800

$ws.Range("H17").Value = 2972.375
$ws.Range("L17").Value = 9762.428400000001
$ws.Range("N17").Value = -10100.4284
$ws.Range("J17").Value = 3254.1428

$ws.Range("H23").Value = 438.0435
$ws.Range("L23").Value = 1297.59999
$ws.Range("N23").Value = -1767.59999
$ws.Range("J23").Value = 432.53333

$ws.Range("H39").Value = 299.9
$ws.Range("M39").Value = -605.6999999999999
$ws.Range("I39").Value = 299.9
$ws.Range("K39").Value = 899.6999999999999

$ws.Range("H55").Value = 7364.9697
$ws.Range("L55").Value = 25302.963
$ws.Range("N55").Value = -25656.963
$ws.Range("J55").Value = 8434.321

$ws.Range("H126").Value = 16169.857
$ws.Range("L126").Value = 79665
$ws.Range("N126").Value = -89545
$ws.Range("M126").Value = -2029
$ws.Range("I126").Value = 2323
$ws.Range("J126").Value = 26555
$ws.Range("K126").Value = 6969

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 10155.412
$ws.Range("M102").Value = -9889.385
$ws.Range("I102").Value = 11511.385
$ws.Range("K102").Value = 11511.385

$ws.Range("H106").Value = 1750000
$ws.Range("L106").Value = 1750000
$ws.Range("N106").Value = -1752524
$ws.Range("J106").Value = 1750000

$ws.Range("H135").Value = 250047630
$ws.Range("L135").Value = 63500
$ws.Range("N135").Value = -73640
$ws.Range("J135").Value = 63500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7618.3335
$ws.Range("L16").Value = 4155
$ws.Range("N16").Value = -4495
$ws.Range("M16").Value = -8530.625
$ws.Range("I16").Value = 8700.625
$ws.Range("J16").Value = 4155
$ws.Range("K16").Value = 8700.625

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H46").Value = 3670540.2
$ws.Range("L46").Value = 5364312.5
$ws.Range("N46").Value = -5364688.5
$ws.Range("M46").Value = -512
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 5364312.5
$ws.Range("K46").Value = 700

$ws.Range("H61").Value = 4452.4243
$ws.Range("L61").Value = 10639
$ws.Range("N61").Value = -11043
$ws.Range("M61").Value = -2875.6296
$ws.Range("I61").Value = 3077.6296
$ws.Range("J61").Value = 10639
$ws.Range("K61").Value = 3077.6296

$ws.Range("H93").Value = 5917.7915
$ws.Range("M93").Value = -5100.9473
$ws.Range("I93").Value = 6348.9473
$ws.Range("K93").Value = 6348.9473

$ws.Range("H100").Value = 3816.2354
$ws.Range("L100").Value = 5165.1113
$ws.Range("N100").Value = -6247.1113
$ws.Range("M100").Value = -1757.75
$ws.Range("I100").Value = 2298.75
$ws.Range("J100").Value = 5165.1113
$ws.Range("K100").Value = 2298.75

$ws.Range("H113").Value = 4452.4243
$ws.Range("L113").Value = 10639
$ws.Range("N113").Value = -14979
$ws.Range("M113").Value = -907.6296000000002
$ws.Range("I113").Value = 3077.6296
$ws.Range("J113").Value = 10639
$ws.Range("K113").Value = 3077.6296

$ws.Range("H136").Value = 3714.3333
$ws.Range("M136").Value = -4941
$ws.Range("I136").Value = 2497
$ws.Range("K136").Value = 7491

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 16382.913
$ws.Range("L107").Value = 302848.008
$ws.Range("N107").Value = -306688.008
$ws.Range("M107").Value = -9173.849999999999
$ws.Range("I107").Value = 3697.95
$ws.Range("J107").Value = 100949.336
$ws.Range("K107").Value = 11093.85

$ws.Range("H122").Value = 3889.0967
$ws.Range("M122").Value = -1744.1578
$ws.Range("I122").Value = 1398.0526
$ws.Range("K122").Value = 4194.1578
